$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextValue "D2" "66.066.21"
$ws.Range("E2").Value = "  -2.23%  "

Set-TextValue "D3" "3.454.30"
$ws.Range("E3").Value = "  -0.67%  "

$ws.Range("E4").Value = "  +0.01%  "

Set-TextValue "D5" "584.91"
$ws.Range("E5").Value = "  -1.11%  "

Set-TextValue "D6" "174.07"
$ws.Range("E6").Value = "  -3.07%  "

$ws.Range("E7").Value = "  +0.02%  "

Set-TextValue "D8" "0.604"
$ws.Range("E8").Value = "  -1.54%  "

Set-TextValue "D9" "3.451.51"
$ws.Range("E9").Value = "  -0.67%  "

$ws.Range("E10").Value = "  -4.65%  "

Set-TextValue "D11" "6.91"
$ws.Range("E11").Value = "  -1.09%  "

Set-TextValue "D12" "0.411"
$ws.Range("E12").Value = "  -3.64%  "

Set-TextValue "D13" "4.055.45"
$ws.Range("E13").Value = "  -0.65%  "

Set-TextValue "D15" "29.05"
$ws.Range("E15").Value = "  -9.79%  "

Set-TextValue "D16" "66.110.35"
$ws.Range("E16").Value = "  -2.15%  "

$ws.Range("E17").Value = "  -2.83%  "

Set-TextValue "D18" "3.451.54"
$ws.Range("E18").Value = "  -0.74%  "

$ws.Range("E19").Value = "  -3.08%  "

Set-TextValue "D20" "13.87"
$ws.Range("E20").Value = "  -1.41%  "

Set-TextValue "D21" "368.72"
$ws.Range("E21").Value = "  -4.52%  "

$ws.Range("E22").Value = "  -2.84%  "

Set-TextValue "D23" "72.74"
$ws.Range("E23").Value = "  +0.84%  "

$ws.Range("E24").Value = "  +0.11%  "

Set-TextValue "D25" "0.538"
$ws.Range("E25").Value = "  +0.38%  "

$ws.Range("E26").Value = "  +0.26%  "

Set-TextValue "D27" "9.77"
$ws.Range("E27").Value = "  -3.05%  "

$ws.Range("E28").Value = "  +0.96%  "

Set-TextValue "D29" "0.999"
$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("E30").Value = "  -2.38%  "

$ws.Range("E31").Value = "  -2.75%  "

Set-TextValue "D32" "5.73"
$ws.Range("E32").Value = "  -5.00%  "

$ws.Range("E34").Value = "  -5.81%  "

Set-TextValue "D35" "7.02"
$ws.Range("E35").Value = "  -3.57%  "

Set-TextValue "D36" "1.54"
$ws.Range("E36").Value = "  -1.61%  "

Set-TextValue "D37" "161.09"
$ws.Range("E37").Value = "  +0.56%  "

Set-TextValue "D38" "28.97"
$ws.Range("E38").Value = "  +4.70%  "

Set-TextValue "D39" "0.882"
$ws.Range("E39").Value = "  -0.48%  "

Set-TextValue "D40" "2.64"
$ws.Range("E40").Value = "  -1.74%  "

$ws.Range("E41").Value = "  -4.42%  "

Set-TextValue "D42" "2.762.50"
$ws.Range("E42").Value = "  +1.37%  "

$ws.Range("E43").Value = "  -1.53%  "

Set-TextValue "D44" "6.46"
$ws.Range("E44").Value = "  -2.74%  "

Set-TextValue "D45" "0.0683"
$ws.Range("E45").Value = "  -3.74%  "

Set-TextValue "D46" "40.20"
$ws.Range("E46").Value = "  -3.14%  "

Set-TextValue "D47" "24.38"
$ws.Range("E47").Value = "  -5.67%  "

Set-TextValue "D48" "0.0291"
$ws.Range("E48").Value = "  -2.07%  "

Set-TextValue "D49" "326.90"
$ws.Range("E49").Value = "  -0.29%  "

Set-TextValue "D50" "6.26"
$ws.Range("E50").Value = "  -0.36%  "

Set-TextValue "D51" "0.102"
$ws.Range("E51").Value = "  -2.70%  "
